# Update countries & provincias Spain
# - Reorders three country-name pairs in the shared-string table (which,
#   because the underlying data rows keep referencing the same string
#   index, manifests as those row labels swapping) and refreshes the
#   "datos actualizados" timestamp.
# - Updates the covid-19 case statistics for a number of countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 10:36"

# --- Country-name swaps (rows keep same position, labels exchange) ----
# Puerto Rico (row 189) <-> Namibia (row 190)
$ws.Range("A189").Value = "Namibia"
$ws.Range("A190").Value = "Puerto Rico"

# Santa Sede (row 208) <-> Islas Turcas y Caicos (row 209)
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"

# Papua Nueva Guinea (row 213) <-> Islas Virgenes Britanicas (row 214)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("A214").Value = "Papua Nueva Guinea"

# --- Updated statistics -------------------------------------------------
# Rusia
$ws.Range("B6").Value = 569063
$ws.Range("C6").Value = 7972
$ws.Range("D6").Value = 324406
$ws.Range("E6").Value = 236816
$ws.Range("G6").Value = 181
$ws.Range("H6").Value = 7841

# Singapur
$ws.Range("B34").Value = 41615
$ws.Range("C34").Value = 142
$ws.Range("E34").Value = 8877

# Filipinas
$ws.Range("B41").Value = 28459
$ws.Range("C41").Value = 660
$ws.Range("D41").Value = 7378
$ws.Range("E41").Value = 19951
$ws.Range("G41").Value = 14
$ws.Range("H41").Value = 1130

# Israel
$ws.Range("B51").Value = 20243
$ws.Range("C51").Value = 207
$ws.Range("D51").Value = 15567
$ws.Range("E51").Value = 4372
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 304

# Kazajistan
$ws.Range("D56").Value = 10195
$ws.Range("E56").Value = 6051

# Estonia
$ws.Range("B103").Value = 1979
$ws.Range("C103").Value = 2
$ws.Range("D103").Value = 1755
$ws.Range("E103").Value = 155

# Eslovaquia
$ws.Range("B114").Value = 1576
$ws.Range("C114").Value = 14
$ws.Range("D114").Value = 1447
$ws.Range("E114").Value = 101

# Estado de Palestina
$ws.Range("B145").Value = 620
$ws.Range("C145").Value = 20
$ws.Range("E145").Value = 202

# Row 189 (now Namibia)
$ws.Range("B189").Value = 40
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 19
$ws.Range("E189").Value = 21
$ws.Range("H189").Value = 0

# Row 190 (now Puerto Rico)
$ws.Range("D190").Value = 1
$ws.Range("E190").Value = 36
$ws.Range("H190").Value = 2

# Row 208 (now Islas Turcas y Caicos)
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Row 209 (now Santa Sede)
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Row 213 (now Islas Virgenes Britanicas)
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

# Row 214 (now Papua Nueva Guinea)
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
